$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5444.222
$ws.Range("I74").Value = 4332.6665
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 4332.6665
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -3396.6665
$ws.Range("N74").Value = -7872
$ws.Range("H75").Value = 75000
$ws.Range("J75").Value = 75000
$ws.Range("L75").Value = 75000
$ws.Range("N75").Value = -76872
$ws.Range("H76").Value = 5143.8887
$ws.Range("I76").Value = 4624.375
$ws.Range("J76").Value = 9300
$ws.Range("K76").Value = 4624.375
$ws.Range("L76").Value = 9300
$ws.Range("M76").Value = -4309.375
$ws.Range("N76").Value = -9930
$ws.Range("H77").Value = 5444.222
$ws.Range("I77").Value = 4332.6665
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 21663.3325
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -16983.3325
$ws.Range("N77").Value = -39360
$ws.Range("H78").Value = 75000
$ws.Range("J78").Value = 75000
$ws.Range("L78").Value = 225000
$ws.Range("N78").Value = -234360
$ws.Range("H79").Value = 5143.8887
$ws.Range("I79").Value = 4624.375
$ws.Range("J79").Value = 9300
$ws.Range("K79").Value = 4624.375
$ws.Range("L79").Value = 9300
$ws.Range("M79").Value = -3532.375
$ws.Range("N79").Value = -11484
$ws.Range("H80").Value = 675.9355
$ws.Range("I80").Value = 522.9286
$ws.Range("J80").Value = 801.94116
$ws.Range("K80").Value = 1568.7858
$ws.Range("L80").Value = 2405.82348
$ws.Range("M80").Value = -570.7857999999999
$ws.Range("N80").Value = -4401.82348
$ws.Range("H83").Value = 675.9355
$ws.Range("I83").Value = 522.9286
$ws.Range("J83").Value = 801.94116
$ws.Range("K83").Value = 4706.3574
$ws.Range("L83").Value = 7217.47044
$ws.Range("M83").Value = 285.6426000000001
$ws.Range("N83").Value = -17201.47044
$ws.Range("H112").Value = 2895.0227
$ws.Range("I112").Value = 3458.3333
$ws.Range("J112").Value = 2875.141
$ws.Range("K112").Value = 10374.9999
$ws.Range("L112").Value = 8625.423000000001
$ws.Range("M112").Value = -9266.999899999999
$ws.Range("N112").Value = -10841.423
$ws.Range("H125").Value = 3977.3333
$ws.Range("J125").Value = 5519.8
$ws.Range("L125").Value = 49678.2
$ws.Range("N125").Value = -54598.2
$ws.Range("H138").Value = 6858.1035
$ws.Range("I138").Value = 3251.9167
$ws.Range("J138").Value = 7798.8477
$ws.Range("K138").Value = 9755.750100000001
$ws.Range("L138").Value = 23396.5431
$ws.Range("M138").Value = -4615.750100000001
$ws.Range("N138").Value = -33676.5431

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 15000
$ws.Range("I39").Value = 15000
$ws.Range("K39").Value = 15000
$ws.Range("M39").Value = -14480
$ws.Range("H61").Value = 47623840
$ws.Range("I61").Value = 71431450
$ws.Range("K61").Value = 71431450
$ws.Range("M61").Value = -71431238
$ws.Range("H74").Value = 111236770
$ws.Range("I74").Value = 111236770
$ws.Range("K74").Value = 111236770
$ws.Range("M74").Value = -111235896
$ws.Range("H77").Value = 111236770
$ws.Range("I77").Value = 111236770
$ws.Range("K77").Value = 556183850
$ws.Range("M77").Value = -556179482
$ws.Range("H132").Value = 66671480
$ws.Range("I132").Value = 4701.9165
$ws.Range("K132").Value = 14105.7495
$ws.Range("M132").Value = -11575.7495
$ws.Range("H136").Value = 47623840
$ws.Range("I136").Value = 71431450
$ws.Range("K136").Value = 214294350
$ws.Range("M136").Value = -214291800

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 1497.8889
$ws.Range("I37").Value = 1662
$ws.Range("K37").Value = 1662
$ws.Range("M37").Value = -1525
$ws.Range("H134").Value = 1591
$ws.Range("I134").Value = 1591
$ws.Range("K134").Value = 4773
$ws.Range("M134").Value = -2238

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 4657.625
$ws.Range("I22").Value = 6249
$ws.Range("J22").Value = 792.8570999999999
$ws.Range("K22").Value = 6249
$ws.Range("L22").Value = 792.8570999999999
$ws.Range("M22").Value = -5899
$ws.Range("N22").Value = -1492.8571
$ws.Range("H31").Value = 21281072
$ws.Range("I31").Value = 3266.9062
$ws.Range("J31").Value = 66673720
$ws.Range("K31").Value = 3266.9062
$ws.Range("L31").Value = 66673720
$ws.Range("M31").Value = -2971.9062
$ws.Range("N31").Value = -66674310
$ws.Range("H34").Value = 21281072
$ws.Range("I34").Value = 3266.9062
$ws.Range("J34").Value = 66673720
$ws.Range("K34").Value = 3266.9062
$ws.Range("L34").Value = 66673720
$ws.Range("M34").Value = -3064.9062
$ws.Range("N34").Value = -66674124
$ws.Range("H58").Value = 1429.2258
$ws.Range("I58").Value = 1393.5333
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 1393.5333
$ws.Range("L58").Value = 2500
$ws.Range("M58").Value = -1190.5333
$ws.Range("N58").Value = -2906
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 2556.5588
$ws.Range("I132").Value = 2058.3
$ws.Range("K132").Value = 6174.900000000001
$ws.Range("M132").Value = -3644.900000000001
$ws.Range("H136").Value = 1429.2258
$ws.Range("I136").Value = 1393.5333
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 4180.5999
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -1630.5999
$ws.Range("N136").Value = -12600
$ws.Range("H141").Value = 336457.34
$ws.Range("J141").Value = 347280.2
$ws.Range("L141").Value = 347280.2
$ws.Range("N141").Value = -357640.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1616.0769
$ws.Range("I5").Value = 1030
$ws.Range("K5").Value = 3090
$ws.Range("M5").Value = -2978
$ws.Range("H128").Value = 116138.75
$ws.Range("I128").Value = 116138.75
$ws.Range("K128").Value = 348416.25
$ws.Range("M128").Value = -343436.25
$ws.Range("H131").Value = 1750.2632
$ws.Range("J131").Value = 1835
$ws.Range("L131").Value = 5505
$ws.Range("N131").Value = -15585
$ws.Range("H132").Value = 3707580.8
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 36000
$ws.Range("M132").Value = -33470
$ws.Range("H133").Value = 12692
$ws.Range("J133").Value = 19959.5
$ws.Range("L133").Value = 59878.5
$ws.Range("N133").Value = -69998.5
$ws.Range("H135").Value = 1616.0769
$ws.Range("I135").Value = 1030
$ws.Range("K135").Value = 9270
$ws.Range("M135").Value = -6735

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 812.4
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H82").Value = 3152.6428
$ws.Range("I82").Value = 1126.8889
$ws.Range("J82").Value = 6799
$ws.Range("K82").Value = 1126.8889
$ws.Range("L82").Value = 6799
$ws.Range("M82").Value = -765.8888999999999
$ws.Range("N82").Value = -7521
$ws.Range("H85").Value = 3152.6428
$ws.Range("I85").Value = 1126.8889
$ws.Range("J85").Value = 6799
$ws.Range("K85").Value = 1126.8889
$ws.Range("L85").Value = 6799
$ws.Range("M85").Value = 121.1111000000001
$ws.Range("N85").Value = -9295
$ws.Range("H124").Value = 46570.6
$ws.Range("J124").Value = 46570.6
$ws.Range("L124").Value = 46570.6
$ws.Range("N124").Value = -56390.6
$ws.Range("H136").Value = 6142.1304
$ws.Range("I136").Value = 5264.65
$ws.Range("K136").Value = 15793.95
$ws.Range("M136").Value = -13243.95

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 14347703
$ws.Range("J135").Value = 14347703
$ws.Range("L135").Value = 14347703
$ws.Range("N135").Value = -14357843
